$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.110.48"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.319.06"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.02"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.76"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "2.679.81"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "2.308.35"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D18").Value = "43.018.90"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.71"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.96"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.34"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.32%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "1.995.21"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0289"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -9.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.59"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.92"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.90"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "2.545.37"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  +0.05%  "
